$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.1716463333333333
$ws.Range("H2").Value = 0.514939
$ws.Range("I2").Value = 0.07772289907851986
$ws.Range("J2").Value = 0.07772289907851984
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.03970866666666666
$ws.Range("N2").Value = 0.119126
$ws.Range("Q2").Value = 0.006815847034888889
$ws.Range("R2").Value = 0.06134262331400001
$ws.Range("S2").Value = 0.07772289907851986
$ws.Range("T2").Value = 0.07772289907851984

# Row 3
$ws.Range("I3").Value = 0.652463224789629
$ws.Range("J3").Value = 0.652463224789629
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.03970866666666666
$ws.Range("N3").Value = 0.119126
$ws.Range("Q3").Value = 0.05721723698911111
$ws.Range("R3").Value = 0.514955132902
$ws.Range("S3").Value = 0.652463224789629
$ws.Range("T3").Value = 0.652463224789629

# Row 4
$ws.Range("G4").Value = 0.5958676666666666
$ws.Range("H4").Value = 1.787603
$ws.Range("I4").Value = 0.2698138761318511
$ws.Range("J4").Value = 0.2698138761318511
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.03970866666666666
$ws.Range("N4").Value = 0.119126
$ws.Range("Q4").Value = 0.0236611105531111
$ws.Range("R4").Value = 0.212949994978
$ws.Range("S4").Value = 0.2698138761318511
$ws.Range("T4").Value = 0.2698138761318511
